$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "WEB -> Internet" firewall rule (row 3) was removed entirely, shifting every
# row below it up by one.
$ws.Rows(3).Delete()

# The old "Internet -> LB" rule (originally row 6, now row 5 after the shift above)
# was also removed, but as a content-clear rather than a row delete, so the blank
# row stays in place (dimension ends at row 15, not 14).
$ws.Range("A5:J5").ClearContents()

# The two remaining rule rows (originally numbered 20 and 30 in column A, for the
# WAS and DB rules) were renumbered to 4 and 5.
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 5

# Restore the user's active selection/cursor to D9 (a single cell, not a range).
$ws.Range("D9").Select()
